# fix: Recreate 10 truly unique Kimi-style sample templates
# Update the title/body text boxes on slides 1 and 2: left-align the
# paragraphs and switch the font from Calibri to Arial.

$ppAlignLeft = 1

$p = $ppt.ActivePresentation

# ---- Slide 1 -----------------------------------------------------------
$s1 = $p.Slides.Item(1)

# "Ocean Blue" title textbox
$tr = $s1.Shapes.Item(4).TextFrame.TextRange
$tr.ParagraphFormat.Alignment = $ppAlignLeft
$tr.Font.Name = "Arial"

# "Calm & Professional" subtitle textbox
$tr = $s1.Shapes.Item(5).TextFrame.TextRange
$tr.ParagraphFormat.Alignment = $ppAlignLeft
$tr.Font.Name = "Arial"

# ---- Slide 2 -----------------------------------------------------------
$s2 = $p.Slides.Item(2)

# "Blue Features" heading textbox
$tr = $s2.Shapes.Item(3).TextFrame.TextRange
$tr.ParagraphFormat.Alignment = $ppAlignLeft
$tr.Font.Name = "Arial"

# Bullet list textbox
$tr = $s2.Shapes.Item(5).TextFrame.TextRange
$tr.ParagraphFormat.Alignment = $ppAlignLeft
$tr.Font.Name = "Arial"
